# Weekly update: insert two new daily price rows (2022-04-08) for Higo
# right after the existing 2021-05-28 row, pushing all subsequent rows
# down by two (old row 8 -> new row 10, ... old row 28 -> new row 30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 8 (formatting is inherited from the
# row above, which already carries the date-number-format style on column D).
$ws.Range("A8:A9").EntireRow.Insert()

# New row 8: Higo, "Primera" quality
$ws.Cells.Item(8,1).Value2  = 6
$ws.Cells.Item(8,2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(8,3).Value2  = "Metropolitana"
$ws.Cells.Item(8,4).Value2  = 44659
$ws.Cells.Item(8,5).Value2  = 13
$ws.Cells.Item(8,6).Value2  = "Fruta"
$ws.Cells.Item(8,7).Value2  = 100101
$ws.Cells.Item(8,8).Value2  = "Berries"
$ws.Cells.Item(8,9).Value2  = 100101006
$ws.Cells.Item(8,10).Value2 = "Higo"
$ws.Cells.Item(8,11).Value2 = "Sin especificar"
$ws.Cells.Item(8,12).Value2 = "Primera"
$ws.Cells.Item(8,13).Value2 = 50
$ws.Cells.Item(8,14).Value2 = 15000
$ws.Cells.Item(8,15).Value2 = 15000
$ws.Cells.Item(8,16).Value2 = 15000
$ws.Cells.Item(8,17).Value2 = "`$/bandeja 7 kilos"
$ws.Cells.Item(8,18).Value2 = "Región Metropolitana"
$ws.Cells.Item(8,19).Value2 = 2143
$ws.Cells.Item(8,20).Value2 = 7

# New row 9: Higo, "Segunda" quality
$ws.Cells.Item(9,1).Value2  = 6
$ws.Cells.Item(9,2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(9,3).Value2  = "Metropolitana"
$ws.Cells.Item(9,4).Value2  = 44659
$ws.Cells.Item(9,5).Value2  = 13
$ws.Cells.Item(9,6).Value2  = "Fruta"
$ws.Cells.Item(9,7).Value2  = 100101
$ws.Cells.Item(9,8).Value2  = "Berries"
$ws.Cells.Item(9,9).Value2  = 100101006
$ws.Cells.Item(9,10).Value2 = "Higo"
$ws.Cells.Item(9,11).Value2 = "Sin especificar"
$ws.Cells.Item(9,12).Value2 = "Segunda"
$ws.Cells.Item(9,13).Value2 = 20
$ws.Cells.Item(9,14).Value2 = 12000
$ws.Cells.Item(9,15).Value2 = 12000
$ws.Cells.Item(9,16).Value2 = 12000
$ws.Cells.Item(9,17).Value2 = "`$/bandeja 7 kilos"
$ws.Cells.Item(9,18).Value2 = "Región Metropolitana"
$ws.Cells.Item(9,19).Value2 = 1714
$ws.Cells.Item(9,20).Value2 = 7
